$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, matching the source
# workbook's inline-string cells (avoids Excel auto-coercing numeric-
# looking strings such as "6.60" or "1.00" into numbers, which would
# drop trailing zeros / change the stored cell type).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '61.034.63'
Set-TextValue $ws.Range('E2') '  +0.93%  '
Set-TextValue $ws.Range('D3') '2.682.02'
Set-TextValue $ws.Range('E3') '  +3.08%  '
Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  -0.14%  '
Set-TextValue $ws.Range('D5') '581.59'
Set-TextValue $ws.Range('E5') '  +1.66%  '
Set-TextValue $ws.Range('D6') '145.62'
Set-TextValue $ws.Range('E6') '  +2.09%  '
Set-TextValue $ws.Range('D7') '0.997'
Set-TextValue $ws.Range('E7') '  +0.07%  '
Set-TextValue $ws.Range('D8') '0.601'
Set-TextValue $ws.Range('E8') '  +0.18%  '
Set-TextValue $ws.Range('D9') '6.60'
Set-TextValue $ws.Range('E9') '  +1.78%  '
Set-TextValue $ws.Range('E10') '  +2.02%  '
Set-TextValue $ws.Range('E11') '  +4.61%  '
Set-TextValue $ws.Range('E12') '  +1.67%  '
Set-TextValue $ws.Range('D13') '3.134.97'
Set-TextValue $ws.Range('E13') '  +2.21%  '
Set-TextValue $ws.Range('D14') '25.74'
Set-TextValue $ws.Range('E14') '  +10.53%  '
Set-TextValue $ws.Range('D15') '61.011.00'
Set-TextValue $ws.Range('E15') '  +0.82%  '
Set-TextValue $ws.Range('E16') '  +2.19%  '
Set-TextValue $ws.Range('D17') '2.669.82'
Set-TextValue $ws.Range('E17') '  +1.86%  '
Set-TextValue $ws.Range('D18') '11.62'
Set-TextValue $ws.Range('E18') '  +2.31%  '
Set-TextValue $ws.Range('D19') '4.75'
Set-TextValue $ws.Range('E19') '  +2.06%  '
Set-TextValue $ws.Range('E20') '  +1.58%  '
Set-TextValue $ws.Range('D21') '6.97'
Set-TextValue $ws.Range('E21') '  -0.28%  '
Set-TextValue $ws.Range('D22') '0.999'
Set-TextValue $ws.Range('E22') '  +0.11%  '
Set-TextValue $ws.Range('D23') '0.534'
Set-TextValue $ws.Range('E23') '  +0.15%  '
Set-TextValue $ws.Range('D24') '64.22'
Set-TextValue $ws.Range('E24') '  +1.50%  '
Set-TextValue $ws.Range('B25') 'Binance-PegBSC-USD'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range('D25') '0.998'
Set-TextValue $ws.Range('E25') '  +0.07%  '
Set-TextValue $ws.Range('B26') 'Kaspa'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D26') '0.162'
Set-TextValue $ws.Range('E26') '  +1.89%  '
Set-TextValue $ws.Range('D27') '8.23'
Set-TextValue $ws.Range('E27') '  +6.35%  '
Set-TextValue $ws.Range('E28') '  +8.33%  '
Set-TextValue $ws.Range('D29') '0.0₃0822'
Set-TextValue $ws.Range('E29') '  +4.38%  '
Set-TextValue $ws.Range('D30') '6.79'
Set-TextValue $ws.Range('E30') '  +6.04%  '
Set-TextValue $ws.Range('D31') '0.998'
Set-TextValue $ws.Range('E31') '  +0.04%  '
Set-TextValue $ws.Range('D32') '167.06'
Set-TextValue $ws.Range('E32') '  +3.54%  '
Set-TextValue $ws.Range('D33') '19.93'
Set-TextValue $ws.Range('E33') '  +2.26%  '
Set-TextValue $ws.Range('E34') '  +9.26%  '
Set-TextValue $ws.Range('E35') '  +5.93%  '
Set-TextValue $ws.Range('E36') '  +8.27%  '
Set-TextValue $ws.Range('D37') '1.66'
Set-TextValue $ws.Range('E37') '  +3.69%  '
Set-TextValue $ws.Range('D38') '328.50'
Set-TextValue $ws.Range('E38') '  +11.11%  '
Set-TextValue $ws.Range('E39') '  +5.00%  '
Set-TextValue $ws.Range('D40') '38.51'
Set-TextValue $ws.Range('E40') '  +1.76%  '
Set-TextValue $ws.Range('D41') '0.881'
Set-TextValue $ws.Range('E41') '  +3.69%  '
Set-TextValue $ws.Range('E42') '  +6.47%  '
Set-TextValue $ws.Range('D43') '20.64'
Set-TextValue $ws.Range('E43') '  +4.57%  '
Set-TextValue $ws.Range('D44') '135.50'
Set-TextValue $ws.Range('E44') '  -1.51%  '
Set-TextValue $ws.Range('E45') '  +1.92%  '
Set-TextValue $ws.Range('D46') '0.0563'
Set-TextValue $ws.Range('E46') '  +3.31%  '
Set-TextValue $ws.Range('D47') '0.618'
Set-TextValue $ws.Range('E47') '  +1.16%  '
Set-TextValue $ws.Range('D48') '20.57'
Set-TextValue $ws.Range('E48') '  +3.69%  '
Set-TextValue $ws.Range('D49') '0.0248'
Set-TextValue $ws.Range('E49') '  +3.09%  '
Set-TextValue $ws.Range('B50') 'FirstDigitalUSD'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D50') '1.00'
Set-TextValue $ws.Range('E50') '  +0.44%  '
Set-TextValue $ws.Range('B51') 'Maker'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D51') '2.148.43'
Set-TextValue $ws.Range('E51') '  +5.95%  '
